# Update "想去人数" (want-to-go count) and "最低票价" (lowest price) figures
# across the workbook's sheets to reflect freshly scraped data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = 50
$ws1.Range("F3").Value = 7495
$ws1.Range("F4").Value = 282
$ws1.Range("F5").Value = 8
$ws1.Range("F6").Value = 448
$ws1.Range("F7").Value = 4088
$ws1.Range("F11").Value = 655
$ws1.Range("F12").Value = 138

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 8

# --- Sheet "全部类型" (All types, aggregate of the above) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = 50
$ws4.Range("F4").Value = 7495
$ws4.Range("F6").Value = 282
$ws4.Range("F7").Value = 8
$ws4.Range("F8").Value = 448
$ws4.Range("F9").Value = 4088
$ws4.Range("F13").Value = 655
$ws4.Range("F14").Value = 8
$ws4.Range("F15").Value = 138
